$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LIST")

# Clear the values in A3:A6 in place (no row shifting). A3:A5 keep their
# existing style, A6 had no style so clearing it leaves nothing behind.
$ws.Range("A3:A6").ClearContents()

# Update the saved selection on the LIST sheet to C9.
$ws.Range("C9").Select()
